$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("종합설계", 650),
    @("종합설계", 0),
    @("종합설계", 410),
    @("종합설계", 290),
    @("종합설계", -85),
    @("종합설계", 90),
    @("종합설계", 3550)
)

$startRow = 247
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
